# Update model labels (3-way rotation among rows 2, 14, 26)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_14_9_0"
$ws.Range("A14").Value = "model_14_9_12"
$ws.Range("A26").Value = "model_14_9_24"

# New metric values (uniform across all model rows 2-26 after retraining)
$values = @(
    [double]"0.999997025781996",
    [double]"0.9990399586905852",
    [double]"0.9999954579608772",
    [double]"0.9999883959452845",
    [double]"0.999994380528663",
    [double]"2.776301953652792e-06",
    [double]"0.0008961564213956752",
    [double]"9.619536303045964e-07",
    [double]"6.685182888398496e-06",
    [double]"3.823568259351546e-06",
    [double]"0.0001004477823987822",
    [double]"0.001666223860606009",
    [double]"1.000007931248011",
    [double]"0.001737158459910895",
    [double]"91.58878149661076",
    [double]"131.8116837172614"
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
